$d = $word.ActiveDocument

$replacements = @(
    @{Old = "2025-03-11 Tuesday"; New = "2025-03-12 Wednesday"},
    @{Old = "42×48=2016"; New = "84×12=1008"},
    @{Old = "22×47=1034"; New = "65×23=1495"},
    @{Old = "33×48=1584"; New = "12×59=708"},
    @{Old = "32×91=2912"; New = "49×81=3969"},
    @{Old = "94×24=2256"; New = "48×58=2784"},
    @{Old = "45×19=855"; New = "29×43=1247"},
    @{Old = "99×65=6435"; New = "34×46=1564"},
    @{Old = "44×43=1892"; New = "95×53=5035"},
    @{Old = "21×87=1827"; New = "53×74=3922"},
    @{Old = "12×95=1140"; New = "47×62=2914"},
    @{Old = "49×59=2891"; New = "33×45=1485"},
    @{Old = "67×49=3283"; New = "30×52=1560"},
    @{Old = "66×14=924"; New = "91×84=7644"},
    @{Old = "50×31=1550"; New = "67×66=4422"},
    @{Old = "99×84=8316"; New = "18×89=1602"},
    @{Old = "26×52=1352"; New = "47×58=2726"},
    @{Old = "29×66=1914"; New = "87×71=6177"},
    @{Old = "49×96=4704"; New = "65×64=4160"},
    @{Old = "30×18=540"; New = "22×96=2112"},
    @{Old = "22×35=770"; New = "81×75=6075"},
    @{Old = "38×93=3534"; New = "25×53=1325"},
    @{Old = "81×83=6723"; New = "36×63=2268"},
    @{Old = "39×38=1482"; New = "25×48=1200"},
    @{Old = "59×60=3540"; New = "42×37=1554"},
    @{Old = "70×40=2800"; New = "16×32=512"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $found = $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $($r.Old)"
    }
}

Write-Output "Done"
